# This edit rotates the values in column A (rows 1-5) of every data sheet
# ("(r, c)" sheets, i.e. all worksheets except the first, empty "Sheet")
# down by two positions (with wrap-around):
#   A1 <- old A4   A2 <- old A5   A3 <- old A1   A4 <- old A2   A5 <- old A3
# It also nudges the "graph score" value in I1 by a tiny floating point
# amount on the handful of sheets where that recomputed value changed.

$wb = $excel.ActiveWorkbook

for ($i = 2; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Read current column-A values with full fidelity (Text/Value lose
    # precision or don't marshal back correctly for some types, Value2
    # preserves the underlying value for both strings and numbers).
    $a1 = $ws.Range("A1").Value2
    $a2 = $ws.Range("A2").Value2
    $a3 = $ws.Range("A3").Value2
    $a4 = $ws.Range("A4").Value2
    $a5 = $ws.Range("A5").Value2

    # Rotate: each value moves down by two rows, wrapping around.
    $ws.Range("A1").Value = $a4
    $ws.Range("A2").Value = $a5
    $ws.Range("A3").Value = $a1
    $ws.Range("A4").Value = $a2
    $ws.Range("A5").Value = $a3

    # Tiny re-computation precision fix for the "graph score" cell.
    $score = $ws.Range("I1").Value2
    if ($score -eq 184251.8833294109) {
        $ws.Range("I1").Value = 184251.8833294106
    } elseif ($score -eq 184232.9662242235) {
        $ws.Range("I1").Value = 184232.9662242233
    }
}
